# Atualização de bases das ligas, do dia: 29-03-2024 às 17:05
#
# The underlying match records (columns B..AC) for several rows were
# re-sorted / corrected. Column A (the 0-based row index) stays put; the
# record payload (id, teams, odds, ...) moves between rows as follows:
#   - rows 61  <-> 62   (simple swap)
#   - rows 183 -> 184 -> 185 -> 183   (3-way rotation)
#   - rows 186 -> 187 -> 188 -> 186   (3-way rotation)
#   - rows 252 <-> 253  (simple swap)
# plus a handful of standalone odds corrections on rows 265, 266 and 268.
#
# NOTE: row payloads are captured into plain (non-function-returned) local
# variables and written back with Range(...).Value2 = <2-D array> so the
# COM SAFEARRAY shape survives the round trip intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- simple swap: rows 61 and 62 -------------------------------------------
$v61 = $ws.Range("B61:AC61").Value2
$v62 = $ws.Range("B62:AC62").Value2
$ws.Range("B61:AC61").Value2 = $v62
$ws.Range("B62:AC62").Value2 = $v61

# --- 3-way rotation: 183 -> 184 -> 185 -> 183 -------------------------------
# new183 = old184, new184 = old185, new185 = old183
$v183 = $ws.Range("B183:AC183").Value2
$v184 = $ws.Range("B184:AC184").Value2
$v185 = $ws.Range("B185:AC185").Value2
$ws.Range("B183:AC183").Value2 = $v184
$ws.Range("B184:AC184").Value2 = $v185
$ws.Range("B185:AC185").Value2 = $v183

# --- 3-way rotation: 186 -> 187 -> 188 -> 186 -------------------------------
# new186 = old187, new187 = old188, new188 = old186
$v186 = $ws.Range("B186:AC186").Value2
$v187 = $ws.Range("B187:AC187").Value2
$v188 = $ws.Range("B188:AC188").Value2
$ws.Range("B186:AC186").Value2 = $v187
$ws.Range("B187:AC187").Value2 = $v188
$ws.Range("B188:AC188").Value2 = $v186

# --- simple swap: rows 252 and 253 ------------------------------------------
$v252 = $ws.Range("B252:AC252").Value2
$v253 = $ws.Range("B253:AC253").Value2
$ws.Range("B252:AC252").Value2 = $v253
$ws.Range("B253:AC253").Value2 = $v252

# --- standalone odds corrections --------------------------------------------
$ws.Range("U265").Value2 = 1.95
$ws.Range("V265").Value2 = 1.9

$ws.Range("U266").Value2 = 1.875
$ws.Range("V266").Value2 = 1.975

$ws.Range("R268").Value2 = 1.825
$ws.Range("S268").Value2 = 2.025
